$d = $word.ActiveDocument

# --- Paragraph 1: replace "FILLER" with the ACC intro sentence ---
$p1 = $d.Paragraphs(1)

$p1run1 = "The Advocacy and Communications Committee"
$p1run2 = " "
$p1run3 = "(ACC) is a group of dedicated senate representatives that are the student body's voice in the General Senate. "
$p1.Range.Text = $p1run1 + $p1run2 + $p1run3

# --- Insert the 3 new paragraphs after paragraph 1, filling in their text ---
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "We plan events that allow us to listen to student concerns, suggestions, and ideas. We then take that feedback from students, voice it to the Senate and staff, and together we make those changes happen! "

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Additionally, we come up with fun interactive ways for students to become informed on what SGA is, and what we can do for them."

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)

$p4run1 = "Contact: "
$p4run2 = "Chair Liana Nazario at lnazario@floridapoly.edu"
$p4.Range.Text = $p4run1 + $p4run2

# --- Now that all text is in place, split runs by nudging+reverting
#     character formatting at the seams. Doing this last (after every
#     InsertParagraphAfter/Text assignment) keeps the "current typing
#     format" bleed from contaminating later paragraphs. ---

$p1Base = $p1.Range.Start
$spaceStart = $p1Base + $p1run1.Length
$spaceEnd = $spaceStart + $p1run2.Length
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Font.Bold = 1
$spaceRange.Font.Bold = 0

$p4Base = $p4.Range.Start
$contactStart = $p4Base
$contactEnd = $p4Base + $p4run1.Length
$contactRange = $d.Range($contactStart, $contactEnd)
$contactRange.Font.Bold = 1
$contactRange.Font.Bold = 0

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
